$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 98 - this shifts existing rows 98-117 down to 99-118
$ws.Rows.Item(98).EntireRow.Insert()

# Populate the newly inserted row 98 with the new weekly price record
$ws.Range("A98").Value = 11
$ws.Range("B98").Value = 'Vega Monumental Concepción'
$ws.Range("C98").Value = 'Bíobío'
$ws.Range("D98").Value = 44889
$ws.Range("E98").Value = 8
$ws.Range("F98").Value = 'Fruta'
$ws.Range("G98").Value = 100103
$ws.Range("H98").Value = 'Frutos de hueso (carozo)'
$ws.Range("I98").Value = 100103001
$ws.Range("J98").Value = 'Cereza'
$ws.Range("K98").Value = 'Lapins'
$ws.Range("L98").Value = 'Primera'
$ws.Range("M98").Value = 150
$ws.Range("N98").Value = 11000
$ws.Range("O98").Value = 12000
$ws.Range("P98").Value = 11333
$ws.Range("Q98").Value = '$/caja 10 kilos'
$ws.Range("R98").Value = 'Región de Ñuble'
$ws.Range("S98").Value = 1133
$ws.Range("T98").Value = 10
